# Auto-update dashboard 2025-07-22 03:55:25
# Update the "Date" column (B2:B6) on Sheet1 from 2025-06-01 to 2025-07-22,
# keeping the value as literal text (matching the original inline-string
# cells) rather than letting Excel auto-convert the date-like text into a
# serial date number.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$rows = 2, 3, 4, 5, 6
foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 2)
    # Leading apostrophe forces Excel to store the value as text instead of
    # reinterpreting "2025-07-22" as a date serial number.
    $cell.Value = "'2025-07-22"
    # Restore the default "Normal" style so the quote-prefix formatting
    # introduced above doesn't leave a visible style change on the cell.
    $cell.Style = "Normal"
}
